$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2:C81").NumberFormat = "@"

$cArr = New-Object 'object[,]' 80,1
$dArr = New-Object 'object[,]' 80,1
$eArr = New-Object 'object[,]' 80,1

$cArr[0,0] = "01/01/2015"
$dArr[0,0] = 920
$eArr[0,0] = 10.05
$cArr[1,0] = "01/01/2015"
$dArr[1,0] = 2318
$eArr[1,0] = 9.44
$cArr[2,0] = "01/01/2015"
$dArr[2,0] = 1888
$eArr[2,0] = -4.89
$cArr[3,0] = "01/01/2015"
$dArr[3,0] = 1849
$eArr[3,0] = 1.15
$cArr[4,0] = "01/01/2015"
$dArr[4,0] = 2012
$eArr[4,0] = 3.34
$cArr[5,0] = "01/01/2015"
$dArr[5,0] = 1502
$eArr[5,0] = -9.08
$cArr[6,0] = "01/01/2015"
$dArr[6,0] = 2941
$eArr[6,0] = -12.63
$cArr[7,0] = "01/01/2015"
$dArr[7,0] = 3954
$eArr[7,0] = -9.56
$cArr[8,0] = "01/01/2016"
$dArr[8,0] = 837
$eArr[8,0] = -9.02
$cArr[9,0] = "01/01/2016"
$dArr[9,0] = 2352
$eArr[9,0] = 1.47
$cArr[10,0] = "01/01/2016"
$dArr[10,0] = 1779
$eArr[10,0] = -5.77
$cArr[11,0] = "01/01/2016"
$dArr[11,0] = 1705
$eArr[11,0] = -7.79
$cArr[12,0] = "01/01/2016"
$dArr[12,0] = 1796
$eArr[12,0] = -10.74
$cArr[13,0] = "01/01/2016"
$dArr[13,0] = 1453
$eArr[13,0] = -3.26
$cArr[14,0] = "01/01/2016"
$dArr[14,0] = 3199
$eArr[14,0] = 8.77
$cArr[15,0] = "01/01/2016"
$dArr[15,0] = 3854
$eArr[15,0] = -2.53
$cArr[16,0] = "01/01/2017"
$dArr[16,0] = 820
$eArr[16,0] = -2.03
$cArr[17,0] = "01/01/2017"
$dArr[17,0] = 2887
$eArr[17,0] = 22.75
$cArr[18,0] = "01/01/2017"
$dArr[18,0] = 1715
$eArr[18,0] = -3.6
$cArr[19,0] = "01/01/2017"
$dArr[19,0] = 1684
$eArr[19,0] = -1.23
$cArr[20,0] = "01/01/2017"
$dArr[20,0] = 1967
$eArr[20,0] = 9.52
$cArr[21,0] = "01/01/2017"
$dArr[21,0] = 1437
$eArr[21,0] = -1.1
$cArr[22,0] = "01/01/2017"
$dArr[22,0] = 3125
$eArr[22,0] = -2.31
$cArr[23,0] = "01/01/2017"
$dArr[23,0] = 4470
$eArr[23,0] = 15.98
$cArr[24,0] = "01/01/2018"
$dArr[24,0] = 935
$eArr[24,0] = 14.02
$cArr[25,0] = "01/01/2018"
$dArr[25,0] = 2110
$eArr[25,0] = -26.91
$cArr[26,0] = "01/01/2018"
$dArr[26,0] = 1838
$eArr[26,0] = 7.17
$cArr[27,0] = "01/01/2018"
$dArr[27,0] = 1761
$eArr[27,0] = 4.57
$cArr[28,0] = "01/01/2018"
$dArr[28,0] = 1765
$eArr[28,0] = -10.27
$cArr[29,0] = "01/01/2018"
$dArr[29,0] = 1306
$eArr[29,0] = -9.119999999999999
$cArr[30,0] = "01/01/2018"
$dArr[30,0] = 3235
$eArr[30,0] = 3.52
$cArr[31,0] = "01/01/2018"
$dArr[31,0] = 3965
$eArr[31,0] = -11.3
$cArr[32,0] = "01/01/2019"
$dArr[32,0] = 777
$eArr[32,0] = -16.9
$cArr[33,0] = "01/01/2019"
$dArr[33,0] = 1884
$eArr[33,0] = -10.71
$cArr[34,0] = "01/01/2019"
$dArr[34,0] = 1597
$eArr[34,0] = -13.11
$cArr[35,0] = "01/01/2019"
$dArr[35,0] = 1814
$eArr[35,0] = 3.01
$cArr[36,0] = "01/01/2019"
$dArr[36,0] = 2038
$eArr[36,0] = 15.47
$cArr[37,0] = "01/01/2019"
$dArr[37,0] = 1641
$eArr[37,0] = 25.65
$cArr[38,0] = "01/01/2019"
$dArr[38,0] = 2675
$eArr[38,0] = -17.31
$cArr[39,0] = "01/01/2019"
$dArr[39,0] = 3750
$eArr[39,0] = -5.42
$cArr[40,0] = "01/01/2020"
$dArr[40,0] = 781
$eArr[40,0] = 0.51
$cArr[41,0] = "01/01/2020"
$dArr[41,0] = 1760
$eArr[41,0] = -6.58
$cArr[42,0] = "01/01/2020"
$dArr[42,0] = 1452
$eArr[42,0] = -9.08
$cArr[43,0] = "01/01/2020"
$dArr[43,0] = 1896
$eArr[43,0] = 4.52
$cArr[44,0] = "01/01/2020"
$dArr[44,0] = 1797
$eArr[44,0] = -11.83
$cArr[45,0] = "01/01/2020"
$dArr[45,0] = 1613
$eArr[45,0] = -1.71
$cArr[46,0] = "01/01/2020"
$dArr[46,0] = 2566
$eArr[46,0] = -4.07
$cArr[47,0] = "01/01/2020"
$dArr[47,0] = 4024
$eArr[47,0] = 7.31
$cArr[48,0] = "01/01/2021"
$dArr[48,0] = 830
$eArr[48,0] = 6.27
$cArr[49,0] = "01/01/2021"
$dArr[49,0] = 1872
$eArr[49,0] = 6.36
$cArr[50,0] = "01/01/2021"
$dArr[50,0] = 1523
$eArr[50,0] = 4.89
$cArr[51,0] = "01/01/2021"
$dArr[51,0] = 1626
$eArr[51,0] = -14.24
$cArr[52,0] = "01/01/2021"
$dArr[52,0] = 1789
$eArr[52,0] = -0.45
$cArr[53,0] = "01/01/2021"
$dArr[53,0] = 1166
$eArr[53,0] = -27.71
$cArr[54,0] = "01/01/2021"
$dArr[54,0] = 2939
$eArr[54,0] = 14.54
$cArr[55,0] = "01/01/2021"
$dArr[55,0] = 3795
$eArr[55,0] = -5.69
$cArr[56,0] = "01/01/2022"
$dArr[56,0] = 842
$eArr[56,0] = 1.45
$cArr[57,0] = "01/01/2022"
$dArr[57,0] = 1696
$eArr[57,0] = -9.4
$cArr[58,0] = "01/01/2022"
$dArr[58,0] = 1315
$eArr[58,0] = -13.66
$cArr[59,0] = "01/01/2022"
$dArr[59,0] = 1500
$eArr[59,0] = -7.75
$cArr[60,0] = "01/01/2022"
$dArr[60,0] = 1536
$eArr[60,0] = -14.14
$cArr[61,0] = "01/01/2022"
$dArr[61,0] = 1304
$eArr[61,0] = 11.84
$cArr[62,0] = "01/01/2022"
$dArr[62,0] = 2230
$eArr[62,0] = -24.12
$cArr[63,0] = "01/01/2022"
$dArr[63,0] = 3863
$eArr[63,0] = 1.79
$cArr[64,0] = "01/01/2023"
$dArr[64,0] = 909
$eArr[64,0] = 7.96
$cArr[65,0] = "01/01/2023"
$dArr[65,0] = 2230
$eArr[65,0] = 31.49
$cArr[66,0] = "01/01/2023"
$dArr[66,0] = 1503
$eArr[66,0] = 14.3
$cArr[67,0] = "01/01/2023"
$dArr[67,0] = 1687
$eArr[67,0] = 12.47
$cArr[68,0] = "01/01/2023"
$dArr[68,0] = 1867
$eArr[68,0] = 21.55
$cArr[69,0] = "01/01/2023"
$dArr[69,0] = 1138
$eArr[69,0] = -12.73
$cArr[70,0] = "01/01/2023"
$dArr[70,0] = 2511
$eArr[70,0] = 12.6
$cArr[71,0] = "01/01/2023"
$dArr[71,0] = 3781
$eArr[71,0] = -2.12
$cArr[72,0] = "01/01/2024"
$dArr[72,0] = 850
$eArr[72,0] = -6.49
$cArr[73,0] = "01/01/2024"
$dArr[73,0] = 2016
$eArr[73,0] = -9.6
$cArr[74,0] = "01/01/2024"
$dArr[74,0] = 1514
$eArr[74,0] = 0.73
$cArr[75,0] = "01/01/2024"
$dArr[75,0] = 1856
$eArr[75,0] = 10.02
$cArr[76,0] = "01/01/2024"
$dArr[76,0] = 2317
$eArr[76,0] = 24.1
$cArr[77,0] = "01/01/2024"
$dArr[77,0] = 1276
$eArr[77,0] = 12.13
$cArr[78,0] = "01/01/2024"
$dArr[78,0] = 2384
$eArr[78,0] = -5.06
$cArr[79,0] = "01/01/2024"
$dArr[79,0] = 3584
$eArr[79,0] = -5.21

$ws.Range("C2:C81").Value = $cArr
$ws.Range("D2:D81").Value = $dArr
$ws.Range("E2:E81").Value = $eArr